$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "The candidate lacks proficiency in NLP, a key requirement for the job. Additionally, the candidate does not have experience with ML prototyping or data preprocessing. However, the candidate has strong experience in computer vision and PyTorch, making them a good fit."
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = "The candidate demonstrates a keen interest in AI and shows willingness to adapt to new cultures. They exhibit good teamwork and problem-solving skills, in addition to showing a desire to contribute to cutting-edge technologies. However, they need to work on improving their adaptability and time management skills."

$ws.Range("C3").Value = "The candidate has strong skills in PyTorch, Python, and Computer Vision, which are essential for the job. However, the candidate lacks NLP experience, a crucial skill for this position."
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = "The candidate has expressed a keen interest in AI, a willingness to adapt to new cultures, and a strong desire to work in Japan. They have also demonstrated good communication and teamwork skills during project challenges, which aligns with the company's requirements for adaptability and problem-solving."
